$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44257
$ws.Range("M2").Value = 100

$ws.Range("D3").Value = 44253
$ws.Range("M3").Value = 160

$ws.Range("D5").Value = 44250
$ws.Range("M5").Value = 200
